# Slide 15, "Content Placeholder 2" shape: merge the last bullet's three
# runs ("...different UDP destination " + "Port2" (bold/blue) + " is used
# for LM") into a single run reading "...different destination UDP is used
# for LM", keeping the first run's (non-bold, tx2 scheme color) formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shape = $s.Shapes.Item(4)
$tr = $shape.TextFrame.TextRange

$paraIndex = 7
$para = $tr.Paragraphs($paraIndex, 1)

$newText = "Does not modify existing TWAMP Light  (which is for DM) procedure as different destination UDP is used for LM"

$oldLeadLen = "Does not modify existing TWAMP Light  (which is for DM) procedure as different UDP destination ".Length

# Overwrite the first run's text in place so it keeps the first run's
# formatting (tx2 scheme color, not bold), producing the full new wording
# followed by the stale leftover text from the old trailing runs.
$lead = $para.Characters(1, $oldLeadLen)
$lead.Text = $newText

# Trim off whatever old text now trails the freshly written text so the
# paragraph collapses back down to a single run.
$tail = $para.Characters($newText.Length + 1, $para.Length - $newText.Length)
$tail.Text = ""
